# Power tolerance analysis update on the "bw factor" sheet:
# adds a brand new 5th asymmetry table (111-220 reflection, rows 64-72)
# and backfills "HHLM spatial chirp" / "HHLM bandwidth" summary rows
# under each of the four pre-existing tables (rows 31-32, 41-42, 51-52,
# 61-62). The new-table block is written first so the new shared strings
# land in the same order the author typed them in (111-220, then the two
# summary labels).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bw factor")
$ws.Activate()

# ---------------------------------------------------------------
# Block 5 (brand new table, "111-220"): header row 64, data rows 65-70,
# summary rows 71-72
# ---------------------------------------------------------------
$ws.Range("A64").Value = "111-220"
$ws.Range("B64").Value = "bragg"
$ws.Range("C64").Value = "asymmetry"
$ws.Range("D64").Value = "b-factor"
$ws.Range("E64").Value = "bandwidth in"
$ws.Range("F64").Value = "bandwidth out"
$ws.Range("G64").Value = "power in"
$ws.Range("H64").Value = "power out"
$ws.Range("I64").Value = "power absorb"

# Row 65: HHLM1
$ws.Range("A65").Value = "HHLM1"
$ws.Range("B65").Value = 12.0368582999333
$ws.Range("C65").Value = 9
$ws.Range("D65").Formula = "=SIN(RADIANS(B65-C65))/SIN(RADIANS(B65+C65))"
$ws.Range("E65").Formula = '=A$11/1000'
$ws.Range("F65").Value = 3.35
$ws.Range("G65").Formula = "=50*0.95^2"
$ws.Range("H65").Formula = "=G65*F65/E65*0.95"
$ws.Range("I65").Formula = "=G65-H65"

# Row 66: HHLM2
$ws.Range("A66").Value = "HHLM2"
$ws.Range("B66").Formula = "=B65"
$ws.Range("C66").Formula = "=-C65"
$ws.Range("D66").Formula = "=SIN(RADIANS(B66-C66))/SIN(RADIANS(B66+C66))"
$ws.Range("E66").Formula = "=F65"
$ws.Range("F66").Value = 3.2509999999999999
$ws.Range("G66").Formula = "=H65"
$ws.Range("I66").Formula = "=G66-H66"

# Row 67: HHLM3
$ws.Range("A67").Value = "HHLM3"
$ws.Range("B67").Value = 19.91
$ws.Range("C67").Value = 16.899999999999999
$ws.Range("F67").Value = 1.891

# Row 68: HHLM4
$ws.Range("A68").Value = "HHLM4"
$ws.Range("B68").Formula = "=B67"
$ws.Range("C68").Formula = "=-C67"
$ws.Range("F68").Value = 1.8280000000000001

# Row 69: C1
$ws.Range("A69").Value = "C1"
$ws.Range("B69").Formula = "=C5"
$ws.Range("C69").Value = 0
$ws.Range("F69").Value = 0.1

# Row 70: C2
$ws.Range("A70").Value = "C2"
$ws.Range("B70").Formula = "=B69"
$ws.Range("C70").Value = -15
$ws.Range("F70").Value = 0.1

# Shared formulas across the new table (set as ranges so Excel folds them
# into <f t="shared"> groups the same way the original table blocks do)
$ws.Range("H66:H70").Formula = "=G66*F66/E66*0.95"
$ws.Range("D67:D70").Formula = "=SIN(RADIANS(B67-C67))/SIN(RADIANS(B67+C67))"
$ws.Range("E67:E70").Formula = "=F66"
$ws.Range("G67:G70").Formula = "=H66"
$ws.Range("I67:I70").Formula = "=G67-H67"

# Row 71-72: HHLM spatial chirp / HHLM bandwidth summary for the new table
$ws.Range("B71").Value = "HHLM spatial chirp (meV/um)"
$ws.Range("C71").Value = [double]"6.77647058825995E-3"

$ws.Range("B72").Value = "HHLM bandwidth (meV)"
$ws.Range("C72").Value = 66.333555879288298

# ---------------------------------------------------------------
# Block 1 (table starting row 24 / "111-440" => HHLM4): summary rows 31-32
# ---------------------------------------------------------------
$ws.Range("B31").Value = "HHLM spatial chirp (meV/um)"
$ws.Range("C31").Value = [double]"-1.1375067569572901E-18"
$ws.Range("C31").NumberFormat = "0.00E+00"

$ws.Range("B32").Value = "HHLM bandwidth (meV)"
$ws.Range("C32").Value = 142.46559719211601

# ---------------------------------------------------------------
# Block 2 (table starting row 34 / "111-333" => HHLM2): summary rows 41-42
# ---------------------------------------------------------------
$ws.Range("B41").Value = "HHLM spatial chirp (meV/um)"
$ws.Range("C41").Value = [double]"-3.6051209103971797E-2"

$ws.Range("B42").Value = "HHLM bandwidth (meV)"
$ws.Range("C42").Value = 213.240960670487

# ---------------------------------------------------------------
# Block 3 (table starting row 44 / "111-331" => HHLM3): summary rows 51-52
# ---------------------------------------------------------------
$ws.Range("B51").Value = "HHLM spatial chirp (meV/um)"
$ws.Range("C51").Value = [double]"-2.5223707918536498E-2"

$ws.Range("B52").Value = "HHLM bandwidth (meV)"
$ws.Range("C52").Value = 223.86280141900701

# ---------------------------------------------------------------
# Block 4 (table starting row 54 / "111-311" => HHLM1): summary rows 61-62
# ---------------------------------------------------------------
$ws.Range("B61").Value = "HHLM spatial chirp (meV/um)"
$ws.Range("C61").Value = [double]"-1.1375067569572901E-18"
$ws.Range("C61").NumberFormat = "0.00E+00"

$ws.Range("B62").Value = "HHLM bandwidth (meV)"
$ws.Range("C62").Value = 223.86280141900701

# ---------------------------------------------------------------
# View state: keep the sheet active, scroll near the new table, and
# select the last-edited cell (B65), matching the author's final view.
# ---------------------------------------------------------------
$ws.Range("A45").Select()
try { $excel.ActiveWindow.ScrollRow = 45 } catch {}
try { $excel.ActiveWindow.TopLeftCell = $ws.Range("A45") } catch {}
$ws.Range("B65").Select()
